$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows with changed values.
# A leading apostrophe on the D-column values forces Excel to store them as
# text (matching the original inline-string/text cell type) instead of
# auto-converting numeric-looking strings (e.g. "1.001") into numbers.
$ws.Range("D2").Value = "'27.252.85"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "'1.902.85"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'306.04"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "'0.5413"
$ws.Range("E7").Value = "  +3.77%  "
$ws.Range("D8").Value = "'0.3804"
$ws.Range("E8").Value = "  +1.13%  "
$ws.Range("D9").Value = "'0.07295"
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("D10").Value = "'22.05"
$ws.Range("E10").Value = "  +4.18%  "
$ws.Range("D11").Value = "'0.9011"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "'0.08188"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "'95.48"
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("D14").Value = "'5.347"
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").Value = "'1.001"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "'14.84"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").Value = "'0.000008648"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").Value = "'1.356.40"
$ws.Range("E18").Value = "  -30.97%  "
$ws.Range("D19").Value = "'0.9986"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").Value = "'27.312.79"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").Value = "'5.045"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").Value = "'10.81"
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("D23").Value = "'6.509"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("D26").Value = "'18.34"
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("D27").Value = "'1.744"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "'116.59"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("D29").Value = "'4.832"
$ws.Range("E29").Value = "  +0.86%  "
$ws.Range("D30").Value = "'4.673"
$ws.Range("E30").Value = "  -3.92%  "
$ws.Range("D31").Value = "'0.09194"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "'0.8275"
$ws.Range("E32").Value = "  +4.43%  "
$ws.Range("D33").Value = "'0.05070"
$ws.Range("E33").Value = "  +0.89%  "
$ws.Range("D34").Value = "'1.224"
$ws.Range("E34").Value = "  +0.69%  "
$ws.Range("D35").Value = "'3.017"
$ws.Range("E35").Value = "  +1.80%  "
$ws.Range("D36").Value = "'3.322"
$ws.Range("E36").Value = "  -3.13%  "
$ws.Range("D37").Value = "'2.695"
$ws.Range("E37").Value = "  +3.12%  "
$ws.Range("D38").Value = "'0.5977"
$ws.Range("E38").Value = "  +4.45%  "
$ws.Range("D39").Value = "'0.01993"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D41").Value = "'9.241"
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("D42").Value = "'6.667"
$ws.Range("E42").Value = "  +1.78%  "
$ws.Range("D43").Value = "'115.87"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("D44").Value = "'0.5140"
$ws.Range("E44").Value = "  +5.84%  "
$ws.Range("D45").Value = "'0.1529"
$ws.Range("E45").Value = "  +1.02%  "
$ws.Range("D46").Value = "'10.27"
$ws.Range("E46").Value = "  +1.51%  "
$ws.Range("D47").Value = "'0.9979"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("D48").Value = "'1.637"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("D49").Value = "'38.02"
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("D50").Value = "'0.06087"
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("D51").Value = "'63.29"
$ws.Range("E51").Value = "  -0.45%  "

# Row 40: only Volume(1h) changes
$ws.Range("E40").Value = "  +0.32%  "

# Rows 24 and 25: Monero and LidoDAOToken swapped positions, with updated price/volume
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").Value = "'148.45"
$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").Value = "'2.306"
$ws.Range("E25").Value = "  -0.27%  "
